{"js": "// Locate the \"Paired Data Module\" paragraph (starts with the bold \"Paired Data Module:\"\n// label) and the following \"Cleaning Module\" paragraph using a text search so we don't\n// depend on fixed paragraph indices.\nconst pairedResults = context.document.body.search(\"Paired Data Modul\", { matchCase: true });\npairedResults.load(\"items\");\nawait context.sync();\n\nif (pairedResults.items.length === 0) {\n  throw new Error('Could not find \"Paired Data Modul\" text in the document.');\n}\n\nconst pairedParagraph = pairedResults.items[0].paragraphs.getFirst();\nconst pairedRange = pairedParagraph.getRange();\n\n// Replace the whole paragraph's OOXML with a version that:\n//  - drops the leading bold \"Paired Data Module:\" + space runs, and\n//  - merges the \"will investigate \" / \"the difference\" / \" in means...\" runs\n//    into a single run (the trailing \" Some other elements...\" run is left as-is).\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">In this </w:t></w:r>' +\n  '<w:r><w:t>module</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  '<w:r><w:t>you</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> will be able to explore paired data, as </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">there is </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">information on 2 runs by the same competitor (unless they DSQ or DNF on the second run). Additionally, </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">you </w:t></w:r>' +\n  '<w:r><w:t>will investigate the difference in means between the 2 runs to determine if racers are on average faster or slower on a specific run.</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> Some other elements of this module are gathering summary statistics, finding a confidence interval, and interpreting your findings.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\npairedRange.insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Remove the whole \"Cleaning Module\" paragraph (the paragraph right after the\n// \"Paired Data Module\" one) entirely, including its paragraph mark.\nconst cleaningResults = context.document.body.search(\"Cleaning Module:\", { matchCase: true });\ncleaningResults.load(\"items\");\nawait context.sync();\n\nif (cleaningResults.items.length > 0) {\n  const cleaningParagraph = cleaningResults.items[0].paragraphs.getFirst();\n  cleaningParagraph.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Paired Data Module\" paragraph by its distinctive leading text.\n$pairedParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Paired Data Modul*\") {\n        $pairedParagraph = $p\n        break\n    }\n}\nif ($pairedParagraph -eq $null) {\n    throw 'Could not find the \"Paired Data Module\" paragraph.'\n}\n\n# Replace that whole paragraph's content with a version that:\n#  - drops the leading bold \"Paired Data Module:\" + space runs, and\n#  - merges the \"will investigate \" / \"the difference\" / \" in means...\" runs\n#    into a single run (the trailing \" Some other elements...\" run is left as-is).\n$newParagraphXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">In this </w:t></w:r><w:r><w:t>module</w:t></w:r><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:r><w:t>you</w:t></w:r><w:r><w:t xml:space=\"preserve\"> will be able to explore paired data, as </w:t></w:r><w:r><w:t xml:space=\"preserve\">there is </w:t></w:r><w:r><w:t xml:space=\"preserve\">information on 2 runs by the same competitor (unless they DSQ or DNF on the second run). Additionally, </w:t></w:r><w:r><w:t xml:space=\"preserve\">you </w:t></w:r><w:r><w:t>will investigate the difference in means between the 2 runs to determine if racers are on average faster or slower on a specific run.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Some other elements of this module are gathering summary statistics, finding a confidence interval, and interpreting your findings.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$pairedParagraph.Range.InsertXML($newParagraphXml)\n\n# Remove the whole \"Cleaning Module\" paragraph entirely, including its paragraph mark.\n$cleaningParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Cleaning Module:*\") {\n        $cleaningParagraph = $p\n        break\n    }\n}\nif ($cleaningParagraph -ne $null) {\n    $cleaningParagraph.Range.Delete()\n}\n"}
